$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20 - this shifts the existing rows 20:48 down to 21:49,
# matching the new weekly data point being added to the top of this data block.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the latest week's record.
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C20").Value = "Arica y Parinacota"
$ws.Range("D20").Value = 44580
$ws.Range("E20").Value = 15
$ws.Range("F20").Value = 100112031
$ws.Range("G20").Value = "Poroto verde"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 1300
$ws.Range("K20").Value = 1300
$ws.Range("L20").Value = 1400
$ws.Range("M20").Value = 1350
$ws.Range("N20").Value = "$/kilo"
$ws.Range("O20").Value = "Región de Arica y Parinacota"
$ws.Range("P20").Value = 1350
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"
